# Applies the cryptos.xlsx price/volume (and a few name/link) updates
# described by the commit diff, cell by cell, using the Excel COM object
# model exposed on $excel / $wb / $ws.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddress, $text) {
    # Force the cell to stay a text cell (matches the original t="inlineStr"/
    # string cells in the workbook) even when the new value looks numeric
    # (e.g. "42.681.50", "1.00", "0.0800"), then restore the cell's original
    # style/number format so no unintended formatting diff is introduced.
    $rng = $ws.Range($cellAddress)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

Set-TextValue "D2" "42.681.50"
Set-TextValue "E2" "  +0.10%  "

Set-TextValue "D3" "2.526.76"
Set-TextValue "E3" "  -0.89%  "

Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.03%  "

Set-TextValue "D5" "314.66"
Set-TextValue "E5" "  +0.63%  "

Set-TextValue "D6" "98.45"
Set-TextValue "E6" "  -2.17%  "

Set-TextValue "E7" "  -1.35%  "

Set-TextValue "E8" "  +0.02%  "

Set-TextValue "D9" "0.516"
Set-TextValue "E9" "  -2.45%  "

Set-TextValue "D10" "35.22"
Set-TextValue "E10" "  -2.78%  "

Set-TextValue "D11" "0.0800"
Set-TextValue "E11" "  -0.57%  "

Set-TextValue "E12" "  +1.13%  "

Set-TextValue "D13" "7.22"
Set-TextValue "E13" "  -2.17%  "

Set-TextValue "D14" "2.919.83"
Set-TextValue "E14" "  -1.02%  "

Set-TextValue "B15" "WrappedEther"
Set-TextValue "C15" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D15" "2.522.68"
Set-TextValue "E15" "  +0.00%  "

Set-TextValue "B16" "Chainlink"
Set-TextValue "C16" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D16" "15.21"
Set-TextValue "E16" "  -6.21%  "

Set-TextValue "D17" "0.810"
Set-TextValue "E17" "  -3.85%  "

Set-TextValue "D18" "42.690.63"
Set-TextValue "E18" "  +0.09%  "

Set-TextValue "D19" "6.58"
Set-TextValue "E19" "  -3.29%  "

Set-TextValue "B20" "ShibaInu"
Set-TextValue "C20" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D20" "0.0₃0940"
Set-TextValue "E20" "  -1.50%  "

Set-TextValue "B21" "InternetComputer(DFINITY)"
Set-TextValue "C21" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D21" "12.14"
Set-TextValue "E21" "  -1.96%  "

Set-TextValue "D22" "68.99"
Set-TextValue "E22" "  -0.26%  "

Set-TextValue "D23" "241.62"
Set-TextValue "E23" "  -0.63%  "

Set-TextValue "D24" "2.85"
Set-TextValue "E24" "  -2.13%  "

Set-TextValue "E25" "  -3.57%  "

Set-TextValue "E26" "  +0.01%  "

Set-TextValue "D27" "25.49"
Set-TextValue "E27" "  -3.63%  "

Set-TextValue "D28" "2.25"
Set-TextValue "E28" "  -4.60%  "

Set-TextValue "D29" "10.00"
Set-TextValue "E29" "  -1.32%  "

Set-TextValue "D30" "37.47"
Set-TextValue "E30" "  -6.45%  "

Set-TextValue "D31" "5.90"
Set-TextValue "E31" "  +3.65%  "

Set-TextValue "D32" "155.27"
Set-TextValue "E32" "  -2.47%  "

Set-TextValue "E33" "  -2.20%  "

Set-TextValue "B34" "WEMIXToken"
Set-TextValue "C34" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D34" "2.64"
Set-TextValue "E34" "  +0.75%  "

Set-TextValue "B35" "Hedera"
Set-TextValue "C35" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D35" "0.0782"
Set-TextValue "E35" "  -2.55%  "

Set-TextValue "D36" "3.13"
Set-TextValue "E36" "  -1.55%  "

Set-TextValue "D37" "1.96"
Set-TextValue "E37" "  -4.59%  "

Set-TextValue "D38" "17.59"
Set-TextValue "E38" "  -2.89%  "

Set-TextValue "E39" "  -2.81%  "

Set-TextValue "E40" "  -1.00%  "

Set-TextValue "D41" "4.23"
Set-TextValue "E41" "  -0.73%  "

Set-TextValue "D42" "21.84"
Set-TextValue "E42" "  +0.04%  "

Set-TextValue "E43" "  -0.09%  "

Set-TextValue "D44" "2.030.92"
Set-TextValue "E44" "  +3.47%  "

Set-TextValue "D45" "0.0296"
Set-TextValue "E45" "  -0.38%  "

Set-TextValue "D46" "3.21"
Set-TextValue "E46" "  -4.43%  "

Set-TextValue "D47" "8.85"
Set-TextValue "E47" "  -1.09%  "

Set-TextValue "D48" "2.770.52"
Set-TextValue "E48" "  -1.16%  "

Set-TextValue "D49" "80.63"
Set-TextValue "E49" "  -0.44%  "

Set-TextValue "E50" "  -2.75%  "

Set-TextValue "D51" "71.85"
Set-TextValue "E51" "  -1.09%  "
